$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2023Q2收支明細")

# Row 13: new cash entry for 偉群's four drop-in game fees
$ws.Range("A13").Value = "2023/06/12"
$ws.Range("C13").Value = 600
$ws.Range("D13").Value = "現金 @偉群 四次臨打費"

# Row 19: update the balance note date
$ws.Range("D19").Value = "2023/06/12 更新餘額"

$ws.Calculate()

$ws.Range("D13").Select() | Out-Null
